$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add module 10 entries: "Data generating process" and "Model residual"
$ws.Range("A65").Value = 10
$ws.Range("B65").Value = "Data generating process"
$ws.Range("C65").Value = "The true underlying causal structure that gives rise to (generates) the data from which you sampled. The data generating process is not known. We use models to try to emulate or approximate the data generating process."

$ws.Range("A66").Value = 10
$ws.Range("B66").Value = "Model residual"
$ws.Range("C66").Value = "The difference between the model predicted value of the outcomee and the observed value. In spatial epidemiology, model residuals can provide clues as to the presence of missing variables that produce spatial patterns"

# Update the selection to match the final active cell
[void]$ws.Range("C66").Select()
